# Commit: ":hammer: change template file part"
# Rename the worksheet from "case_category" to "part" to match the
# new template purpose (Template_Upload_Part.xlsx), and update the
# active cell selection left behind by the editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet "case_category" -> "part"
$ws.Name = "part"

# Move the active selection to B8 (was C5)
$ws.Range("B8").Select()
